# Add new power plant types to the Electricity Source subscript
# (BDSBaPCF - "Boolean Do Suppliers Bid at Peak Capacity Factors")
# Ref: issues #280 and #99

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BDSBaPCF")

# New plant types, appended below the existing data (which currently ends at row 18)
$newPlants = @(
    "hard coal w CCS",
    "natural gas combined cycle w CCS",
    "biomass w CCS",
    "lignite w CCS",
    "small modular reactor",
    "hydrogen"
)

$startRow = 19
for ($idx = 0; $idx -lt $newPlants.Length; $idx++) {
    $r = $startRow + $idx

    # Column A: plant/source name
    $ws.Cells.Item($r, 1).Value = $newPlants[$idx]

    # Column B: Boolean flag (default to 0 / False), with the same "0" integer
    # number format and light accent fill used to highlight the newly added rows
    $cellB = $ws.Cells.Item($r, 2)
    $cellB.Value = 0
    $cellB.NumberFormat = "0"
    $cellB.Interior.ThemeColor = 8
    $cellB.Interior.TintAndShade = 0.79998168889431442
}

# Leave the newly added rows selected, matching the state the sheet was left in
[void]$ws.Range("A19:A24").Select()
